$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false
$wsStage = $wb.Worksheets.Item("StageMapping")
$wsStage.Delete()
